$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @("D2", "243.66"),
    @("D3", "22.98"),
    @("D4", "5.387"),
    @("D5", "0.05944"),
    @("D7", "6.503"),
    @("D8", "0.8107"),
    @("D9", "0.9287"),
    @("B10", "One"),
    @("C10", "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"),
    @("D10", "0.01116"),
    @("E10", "9OneONEBestin24h"),
    @("B11", "WazirX"),
    @("C11", "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"),
    @("D11", "0.1432"),
    @("E11", "10WazirXWRX"),
    @("B12", "MandalaExchangeToken"),
    @("C12", "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"),
    @("D12", "0.07406"),
    @("E12", "11MandalaExchangeTokenMDX"),
    @("B13", "LiechtensteinCryptoassetsExchange"),
    @("C13", "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"),
    @("D13", "0.03247"),
    @("E13", "12LiechtensteinCryptoassetsExchangeLCX"),
    @("B14", "BitrueCoin"),
    @("C14", "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"),
    @("D14", "0.03077"),
    @("E14", "13BitrueCoinBTR"),
    @("B15", "BitMartToken"),
    @("C15", "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"),
    @("D15", "0.09357"),
    @("E15", "14BitMartTokenBMX"),
    @("B16", "MCDex"),
    @("C16", "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"),
    @("D16", "3.864"),
    @("E16", "15MCDexMCB"),
    @("B17", "BitForexToken"),
    @("C17", "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"),
    @("D17", "0.001571"),
    @("E17", "16BitForexTokenBF"),
    @("B18", "CoinExToken"),
    @("C18", "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"),
    @("D18", "0.04701"),
    @("E18", "17CoinExTokenCET"),
    @("B19", "TigerCash"),
    @("C19", "https://coinranking.com/coin/6hIn06L2+tigercash-tch"),
    @("D19", "0.005947"),
    @("E19", "18TigerCashTCH"),
    @("B20", "BitKan"),
    @("C20", "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"),
    @("D20", "0.001256"),
    @("E20", "19BitKanKAN"),
    @("B21", "HotbitToken"),
    @("C21", "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"),
    @("D21", "0.004790"),
    @("E21", "20HotbitTokenHTB"),
    @("B22", "NitroEx"),
    @("C22", "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"),
    @("D22", "0.00007999"),
    @("E22", "21NitroExNTX"),
    @("B23", "LEO"),
    @("C23", "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"),
    @("D23", "3.564"),
    @("E23", "22LEOLEO"),
    @("B24", "BTSEToken"),
    @("C24", "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"),
    @("D24", "2.133"),
    @("E24", "23BTSETokenBTSE"),
    @("D26", "0.1331"),
    @("D27", "0.0002339"),
    @("D40", "0.03923"),
    @("D41", "0.006370"),
    @("D42", "0.1074"),
    @("D44", "0.008202"),
    @("D45", "0.00005183"),
    @("D47", "0.6779"),
    @("D48", "0.002066")
)

foreach ($pair in $changes) {
    $ref = $pair[0]
    $val = $pair[1]
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Write-Output ("Applied " + $changes.Count + " cell updates")
